$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 523, shifting existing rows 523-592 down to 524-593
$ws.Rows.Item(523).Insert()

# Populate the newly inserted row 523 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R are identical to the surrounding rows (constants for this dataset).
$ws.Cells.Item(523, 1).Value = 5
$ws.Cells.Item(523, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(523, 3).Value = "Maule"
$ws.Cells.Item(523, 4).Value = 45124
$ws.Cells.Item(523, 5).Value = 7
$ws.Cells.Item(523, 6).Value = 100112032
$ws.Cells.Item(523, 7).Value = "Zapallo italiano"
$ws.Cells.Item(523, 8).Value = "Sin especificar"
$ws.Cells.Item(523, 9).Value = "Primera"
$ws.Cells.Item(523, 10).Value = 400
$ws.Cells.Item(523, 11).Value = 14000
$ws.Cells.Item(523, 12).Value = 14000
$ws.Cells.Item(523, 13).Value = 14000
$ws.Cells.Item(523, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(523, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(523, 16).Value = 280
$ws.Cells.Item(523, 17).Value = 50
$ws.Cells.Item(523, 18).Value = "Hortaliza"
